# Auto-generated Excel COM-interop script
# Refreshes ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets with updated
# Universalis market-price data (currentAveragePrice / NQ / HQ, leve prices,
# and the derived profit columns) as pulled by the scheduled data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1506.2821
$ws.Range("I15").Value = 1506.2821
$ws.Range("K15").Value = 4518.846299999999
$ws.Range("M15").Value = -4349.846299999999

$ws.Range("H33").Value = 1339.375
$ws.Range("I33").Value = 2532.5
$ws.Range("J33").Value = 146.25
$ws.Range("K33").Value = 2532.5
$ws.Range("L33").Value = 146.25
$ws.Range("M33").Value = -2303.5
$ws.Range("N33").Value = -604.25

$ws.Range("H62").Value = 25004536
$ws.Range("I62").Value = 62502000
$ws.Range("K62").Value = 62502000
$ws.Range("M62").Value = -62501376

$ws.Range("H65").Value = 25004536
$ws.Range("I65").Value = 62502000
$ws.Range("K65").Value = 312510000
$ws.Range("M65").Value = -312506880

$ws.Range("H86").Value = 7524510.5
$ws.Range("I86").Value = 7499.75
$ws.Range("K86").Value = 7499.75
$ws.Range("M86").Value = -6376.75

$ws.Range("H89").Value = 7524510.5
$ws.Range("I89").Value = 7499.75
$ws.Range("K89").Value = 37498.75
$ws.Range("M89").Value = -31882.75

$ws.Range("H125").Value = 2933.125
$ws.Range("I125").Value = 1448.125
$ws.Range("J125").Value = 4418.125
$ws.Range("K125").Value = 13033.125
$ws.Range("L125").Value = 39763.125
$ws.Range("M125").Value = -10573.125
$ws.Range("N125").Value = -44683.125

$ws.Range("H129").Value = 1530.56
$ws.Range("I129").Value = 620.7143
$ws.Range("K129").Value = 1862.1429
$ws.Range("M129").Value = 3137.8571

$ws.Range("H132").Value = 6186.5
$ws.Range("I132").Value = 8394.166999999999
$ws.Range("K132").Value = 25182.501
$ws.Range("M132").Value = -22652.501

$ws.Range("H135").Value = 953777.75
$ws.Range("I135").Value = 1177696
$ws.Range("J135").Value = 2125
$ws.Range("K135").Value = 10599264
$ws.Range("L135").Value = 19125
$ws.Range("M135").Value = -10596729
$ws.Range("N135").Value = -24195

$ws.Range("H137").Value = 1006874.6
$ws.Range("I137").Value = 1113792.9
$ws.Range("J137").Value = 919396
$ws.Range("K137").Value = 3341378.7
$ws.Range("L137").Value = 2758188
$ws.Range("M137").Value = -3338828.7
$ws.Range("N137").Value = -2763288

$ws.Range("H141").Value = 2011.6428
$ws.Range("I141").Value = 2081.28
$ws.Range("K141").Value = 6243.84
$ws.Range("M141").Value = -1063.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4407.768
$ws.Range("I32").Value = 3940.836
$ws.Range("K32").Value = 3940.836
$ws.Range("M32").Value = -3653.836

$ws.Range("H45").Value = 2125.5652
$ws.Range("I45").Value = 2209.4666
$ws.Range("K45").Value = 2209.4666
$ws.Range("M45").Value = -1832.4666

$ws.Range("H102").Value = 1472.6666
$ws.Range("I102").Value = 1486.5714
$ws.Range("J102").Value = 1424
$ws.Range("K102").Value = 1486.5714
$ws.Range("L102").Value = 1424
$ws.Range("M102").Value = 135.4286
$ws.Range("N102").Value = -4668

$ws.Range("H132").Value = 1972.3265
$ws.Range("I132").Value = 1555.238
$ws.Range("J132").Value = 4474.857
$ws.Range("K132").Value = 4665.714
$ws.Range("L132").Value = 13424.571
$ws.Range("M132").Value = -2135.714
$ws.Range("N132").Value = -18484.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 437625.22
$ws.Range("I107").Value = 2888.0527
$ws.Range("K107").Value = 2888.0527
$ws.Range("M107").Value = -968.0527000000002

$ws.Range("H134").Value = 80976.62
$ws.Range("I134").Value = 3269.6
$ws.Range("J134").Value = 340000
$ws.Range("K134").Value = 9808.799999999999
$ws.Range("L134").Value = 1020000
$ws.Range("M134").Value = -7273.799999999999
$ws.Range("N134").Value = -1025070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 800000
$ws.Range("J9").Value = 800000
$ws.Range("L9").Value = 800000
$ws.Range("N9").Value = -800336

$ws.Range("H31").Value = 53763.668
$ws.Range("I31").Value = 2595.2727
$ws.Range("J31").Value = 110048.9
$ws.Range("K31").Value = 2595.2727
$ws.Range("L31").Value = 110048.9
$ws.Range("M31").Value = -2300.2727
$ws.Range("N31").Value = -110638.9

$ws.Range("H34").Value = 53763.668
$ws.Range("I34").Value = 2595.2727
$ws.Range("J34").Value = 110048.9
$ws.Range("K34").Value = 2595.2727
$ws.Range("L34").Value = 110048.9
$ws.Range("M34").Value = -2393.2727
$ws.Range("N34").Value = -110452.9

$ws.Range("H58").Value = 506266.25
$ws.Range("I58").Value = 838627.2
$ws.Range("K58").Value = 838627.2
$ws.Range("M58").Value = -838424.2

$ws.Range("H123").Value = 61666.668
$ws.Range("J123").Value = 61666.668
$ws.Range("L123").Value = 61666.668
$ws.Range("N123").Value = -71466.66800000001

$ws.Range("H134").Value = 669850.4
$ws.Range("I134").Value = 419351.25
$ws.Range("K134").Value = 1258053.75
$ws.Range("M134").Value = -1255518.75

$ws.Range("H136").Value = 506266.25
$ws.Range("I136").Value = 838627.2
$ws.Range("K136").Value = 2515881.6
$ws.Range("M136").Value = -2513331.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 60
$ws.Range("I38").Value = 60
$ws.Range("K38").Value = 180
$ws.Range("M38").Value = 167

$ws.Range("H80").Value = 4451.5
$ws.Range("J80").Value = 4451.5
$ws.Range("L80").Value = 13354.5
$ws.Range("N80").Value = -15226.5

$ws.Range("H83").Value = 4451.5
$ws.Range("J83").Value = 4451.5
$ws.Range("L83").Value = 40063.5
$ws.Range("N83").Value = -49423.5

$ws.Range("H98").Value = 2780.6
$ws.Range("I98").Value = 2985.6
$ws.Range("J98").Value = 2678.1
$ws.Range("K98").Value = 8956.799999999999
$ws.Range("L98").Value = 8034.299999999999
$ws.Range("M98").Value = -7458.799999999999
$ws.Range("N98").Value = -11030.3

$ws.Range("H129").Value = 1377.091
$ws.Range("I129").Value = 483
$ws.Range("J129").Value = 2122.1667
$ws.Range("K129").Value = 1449
$ws.Range("L129").Value = 6366.500100000001
$ws.Range("M129").Value = 3551
$ws.Range("N129").Value = -16366.5001

$ws.Range("H131").Value = 12824487
$ws.Range("I131").Value = 25641830
$ws.Range("J131").Value = 7143.077
$ws.Range("K131").Value = 76925490
$ws.Range("L131").Value = 21429.231
$ws.Range("M131").Value = -76920450
$ws.Range("N131").Value = -31509.231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2085.5
$ws.Range("I102").Value = 1733.762
$ws.Range("K102").Value = 1733.762
$ws.Range("M102").Value = -111.7619999999999

$ws.Range("H107").Value = 1010.68
$ws.Range("I107").Value = 1062.0769
$ws.Range("J107").Value = 955
$ws.Range("K107").Value = 1062.0769
$ws.Range("L107").Value = 955
$ws.Range("M107").Value = 857.9231
$ws.Range("N107").Value = -4795

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H122").Value = 1848.7174
$ws.Range("I122").Value = 1236.1613
$ws.Range("K122").Value = 3708.4839
$ws.Range("M122").Value = -1258.4839

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 459604.53
$ws.Range("I7").Value = 5555.8184
$ws.Range("K7").Value = 5555.8184
$ws.Range("M7").Value = -5443.8184

$ws.Range("H61").Value = 4377.8213
$ws.Range("I61").Value = 3889.45
$ws.Range("J61").Value = 5598.75
$ws.Range("K61").Value = 3889.45
$ws.Range("L61").Value = 5598.75
$ws.Range("M61").Value = -3687.45
$ws.Range("N61").Value = -6002.75

$ws.Range("H113").Value = 4377.8213
$ws.Range("I113").Value = 3889.45
$ws.Range("J113").Value = 5598.75
$ws.Range("K113").Value = 3889.45
$ws.Range("L113").Value = 5598.75
$ws.Range("M113").Value = -1719.45
$ws.Range("N113").Value = -9938.75

$ws.Range("H122").Value = 1055378.6
$ws.Range("J122").Value = 1253750.6
$ws.Range("L122").Value = 3761251.8
$ws.Range("N122").Value = -3766151.8

$ws.Range("H126").Value = 459604.53
$ws.Range("I126").Value = 5555.8184
$ws.Range("K126").Value = 16667.4552
$ws.Range("M126").Value = -14197.4552

$ws.Range("H132").Value = 4866.8887
$ws.Range("I132").Value = 4002
$ws.Range("K132").Value = 12006
$ws.Range("M132").Value = -9476

$ws.Range("H136").Value = 744597.25
$ws.Range("I136").Value = 912209.9
$ws.Range("K136").Value = 2736629.7
$ws.Range("M136").Value = -2734079.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7749.3335
$ws.Range("I62").Value = 9997.666999999999
$ws.Range("J62").Value = 5501
$ws.Range("K62").Value = 9997.666999999999
$ws.Range("L62").Value = 5501
$ws.Range("M62").Value = -9373.666999999999
$ws.Range("N62").Value = -6749

$ws.Range("H65").Value = 7749.3335
$ws.Range("I65").Value = 9997.666999999999
$ws.Range("J65").Value = 5501
$ws.Range("K65").Value = 49988.335
$ws.Range("L65").Value = 27505
$ws.Range("M65").Value = -46868.335
$ws.Range("N65").Value = -33745

$ws.Range("H132").Value = 119770.555
$ws.Range("I132").Value = 13183.75
$ws.Range("J132").Value = 205040
$ws.Range("K132").Value = 39551.25
$ws.Range("L132").Value = 615120
$ws.Range("M132").Value = -37021.25
$ws.Range("N132").Value = -620180

$ws.Range("H136").Value = 11155692
$ws.Range("I136").Value = 18093148
$ws.Range("K136").Value = 54279444
$ws.Range("M136").Value = -54276894
